$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("46").Insert()
Write-Output "done"
